$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "PlanificacionAulas" sheet as a copy of
#    "PlanificacionHorarios" (same layout/formatting), positioned
#    right after it.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PlanificacionHorarios")
$ws2.Copy($null, $ws2)
$ws3 = $wb.ActiveSheet
$ws3.Name = "PlanificacionAulas"

# Remove the two stray formatted cells below the table (rows 16-19
# on PlanificacionHorarios are not present on the new sheet).
$ws3.Rows("16:19").Delete()

# New title text for the "Aulas" (classrooms) planning sheet.
$ws3.Range("B1").Value = "PLANIFICACION AULA XXXX DPTO. DE  INGENIERIA DE SISTEMA PERIODO 20xx-x"

# Column widths on the new sheet (no override on column C; D:I a bit
# narrower than on PlanificacionHorarios).
$ws3.Columns("C").ColumnWidth = $ws2.Columns("B").ColumnWidth
$ws3.Columns("D:I").ColumnWidth = 20

# ------------------------------------------------------------------
# 2. Update selections / active tab: PlanificacionHorarios keeps the
#    selection on B9:C9 (no longer the tab shown on open), while the
#    new PlanificacionAulas sheet becomes the active tab with E6
#    selected.
# ------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("B9:C9").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("E6").Select() | Out-Null
